# Weekly update: insert a new price-record row for "Bruselas (repollito)"
# at Vega Modelo de Temuco, pushing the existing rows 127-137 down to
# 128-138 and adding a brand new entry at row 127.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 127 - this shifts rows 127..137 down to 128..138
# (carrying over formatting, e.g. the date style on column D).
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with the new weekly record.
$ws.Range("A127").Value = 10
$ws.Range("B127").Value = "Vega Modelo de Temuco"
$ws.Range("C127").Value = "La Araucanía"
$ws.Range("D127").Value = 44826
$ws.Range("E127").Value = 9
$ws.Range("F127").Value = 100112035
$ws.Range("G127").Value = "Bruselas (repollito)"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 80
$ws.Range("K127").Value = 24000
$ws.Range("L127").Value = 24000
$ws.Range("M127").Value = 24000
$ws.Range("N127").Value = "$/malla 10 kilos"
$ws.Range("O127").Value = "Región Metropolitana"
$ws.Range("P127").Value = 2400
$ws.Range("Q127").Value = 10
$ws.Range("R127").Value = "Hortaliza"
